# Update gh-pages output (HeFei expo sheet data refresh), applied to both
# the "展览" and "全部类型" worksheets (their data is identical in this workbook).

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Simple "want-to-go" counter (column F) bumps on existing rows ---
    $ws.Range("F3").Value  = 7308
    $ws.Range("F4").Value  = 5595
    $ws.Range("F6").Value  = 174
    $ws.Range("F9").Value  = 109
    $ws.Range("F10").Value = 87
    $ws.Range("F11").Value = 111
    $ws.Range("F12").Value = 204
    $ws.Range("F13").Value = 53
    $ws.Range("F14").Value = 653
    $ws.Range("F15").Value = 351
    $ws.Range("F17").Value = 12

    # --- Insert a brand-new event row before the old row 18, pushing the
    #     old rows 18/19 down to 19/20. Copy row 17's formatting down so the
    #     new row's cells (notably the bordered/centered index column A)
    #     come out identical to the existing rows. ---
    $ws.Rows.Item(18).Insert()
    $ws.Range("A17").Copy($ws.Range("A18"))

    $ws.Range("A18").Value = 17
    # Force text storage for the date-looking string (otherwise the host
    # auto-coerces "2024-05-05" into a date serial number).
    $ws.Range("B18").Value = "'2024-05-05"
    $ws.Range("C18").Value = "合肥·HF动漫展"
    $ws.Range("D18").Value = "文忠路1865号 赫拉诺言艺术中心"
    $ws.Range("E18").Value = "2024.05.05 10:00-05.05 16:00"
    $ws.Range("F18").Value = 0
    $ws.Range("G18").Value = 40
    $ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=83162"
    $ws.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202403/d0StKnDC1710903045879.jpeg"

    # --- The row that used to be r19 (index 18, "运动番only-群青日和") is now
    #     r20; bump its index number and "want-to-go" count. ---
    $ws.Range("A20").Value = 19
    $ws.Range("F20").Value = 47
}
